$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns AD/AE/AF for Wins/Losses/Ties, matching the
# existing header formatting (bold, bordered, centered) used by the other
# header cells such as AC1.
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Team record is the same for every player row (2-44): 92 wins, 70
# losses, 0 ties.
for ($row = 2; $row -le 44; $row++) {
    $ws.Cells.Item($row, 30).Value = 92
    $ws.Cells.Item($row, 31).Value = 70
    $ws.Cells.Item($row, 32).Value = 0
}
